# New Orleans shard workbook rework:
#  - the two worksheets swap roles: the sheet that used to be "hotel_info"
#    (1st tab / rId1) becomes "review_info" (header row only, no data),
#    and the sheet that used to be "review_info" (2nd tab / rId2) becomes
#    "hotel_info" (header row + one data row, with a new "State" column).

$wb = $excel.ActiveWorkbook

# Sheets in tab order - index 1 is rId1 (was "hotel_info"), index 2 is rId2 (was "review_info")
$wsFirst  = $wb.Worksheets.Item(1)
$wsSecond = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1st tab (rId1): wipe the old hotel_info content, replace with the
# review_info header row (no data rows beneath it).
# ---------------------------------------------------------------------
$wsFirst.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $wsFirst.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# ---------------------------------------------------------------------
# 2nd tab (rId2): wipe the old review_info header row, replace with the
# hotel_info header row (now including a "State" column) plus the single
# hotel data row.
# ---------------------------------------------------------------------
$wsSecond.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $wsSecond.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$wsSecond.Cells.Item(2, 1).Value = 41264
$wsSecond.Cells.Item(2, 2).Value = "Extended Stay America New Orleans Kenner"
$wsSecond.Cells.Item(2, 3).Value = "Louisiana"
$wsSecond.Cells.Item(2, 4).Value = "Kenner"
$wsSecond.Cells.Item(2, 5).Value = 70062
$wsSecond.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g40247-d225930-Reviews-Extended_Stay_America_New_Orleans_Airport-Kenner_Louisiana.html"
$wsSecond.Cells.Item(2, 7).Value = "Extended Stay America - New Orleans - Airport"
# These look numeric but must stay text (leading apostrophe keeps them as strings).
$wsSecond.Cells.Item(2, 8).Value = "'126"
$wsSecond.Cells.Item(2, 9).Value = "'13"
$wsSecond.Cells.Item(2, 10).Value = "'127"

# ---------------------------------------------------------------------
# Rename the tabs so their (unchanged) rIds now carry the swapped names:
# rId1 -> "review_info", rId2 -> "hotel_info". Go through a scratch name
# first so the intermediate state never collides with the other tab's
# current name.
# ---------------------------------------------------------------------
$wsFirst.Name = "__tmp_swap__"
$wsSecond.Name = "hotel_info"
$wsFirst.Name = "review_info"
